$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates (PacMan map tile edits) ---
$ws.Range("T2").Value = 1
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 2
$ws.Range("T3").Value = 1
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 2
$ws.Range("T4").Value = 1
$ws.Range("R5").Value = 2
$ws.Range("T5").Value = 1
$ws.Range("P6").Value = 2
$ws.Range("R6").Value = 3
$ws.Range("Q7").Value = 1
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = 1
$ws.Range("R8").Value = 2
$ws.Range("S8").Value = 2
$ws.Range("T8").Value = 1
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 1
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 1
$ws.Range("R10").Value = 1
$ws.Range("T10").Value = 1
$ws.Range("N12").Value = 2
$ws.Range("T12").Value = 1
$ws.Range("S13").Value = 2
$ws.Range("T13").Value = 1
$ws.Range("Q14").Value = 2
$ws.Range("T14").Value = 1
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = 1
$ws.Range("S15").Value = 2
$ws.Range("T15").Value = 1
$ws.Range("R16").Value = 1
$ws.Range("S16").Value = 2
$ws.Range("T16").Value = 1
$ws.Range("S17").Value = 2
$ws.Range("T17").Value = 1
$ws.Range("P18").Value = 3
$ws.Range("S18").Value = 2
$ws.Range("T18").Value = 1
$ws.Range("T19").Value = 1
$ws.Range("B20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("F20").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("J20").Value = 1
$ws.Range("L20").Value = 1
$ws.Range("T20").Value = 1

# --- Rows 21 and 22 are emptied and hidden (map shortened to 20 playable rows) ---
$ws.Range("A21:U22").Clear()
$ws.Rows.Item(21).Hidden = $true
$ws.Rows.Item(22).Hidden = $true

# --- Column U (21) becomes hidden, matching the hidden helper columns V:X ---
$ws.Columns.Item(21).Hidden = $true

# --- Update the saved selection on the sheet ---
$ws.Range("A1:T20").Select()
